$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1): formula, formula2, date_time ---
$ws.Range("E1").Value = "formula"
$ws.Range("E1").Font.ThemeColor = 1

$ws.Range("F1").Value = "formula2"
$ws.Range("F1").Font.ThemeColor = 1

$ws.Range("G1").Value = "date_time"
$ws.Range("G1").Font.ThemeColor = 1

# --- Column E: shared formula A+B ---
$ws.Range("E2:E5").Formula = "=A2+B2"
$ws.Range("E2:E5").Font.ThemeColor = 1

# --- Column F: shared formula E/(A-1), formatted with 3 decimals ---
$ws.Range("F2:F5").Formula = "=E2/(A2-1)"
$ws.Range("F2:F5").Font.ThemeColor = 1
$ws.Range("F2:F5").NumberFormat = "0.000"

# --- Column G: date + time values ---
$ws.Range("G2").Value = 44197.54194444444
$ws.Range("G3").Value = 30371.041944444445
$ws.Range("G4").Value = 44217.0
$ws.Range("G5").Value = 44200.54194444444

$ws.Range("G2:G5").Font.ThemeColor = 1
$ws.Range("G2").NumberFormat = "yyyy. mm. dd h:mm:ss"
$ws.Range("G3").NumberFormat = "yyyy. mm. dd h:mm:ss"
$ws.Range("G4").NumberFormat = "yyyy. m. d h:mm:ss"
$ws.Range("G5").NumberFormat = "yyyy. mm. dd h:mm:ss"

# --- Column widths for F and G ---
$ws.Range("F1:G1").ColumnWidth = 17.25
